$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 ---
$ws.Cells.Item(1,1).Value = 1
$ws.Cells.Item(1,2).Value = "prueba2"
$ws.Cells.Item(1,3).Style = "Normal"
$ws.Cells.Item(1,3).Value = "17-09-2024"
$ws.Cells.Item(1,4).Style = "Normal"
$ws.Cells.Item(1,4).Value = "25-02-2025"
$ws.Cells.Item(1,5).Value = "❌"

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = 2
$ws.Cells.Item(2,2).Value = "Prueba8"
$ws.Cells.Item(2,3).Style = "Normal"
$ws.Cells.Item(2,3).Value = "18-09-2024"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,4).Value = "23-10-2024"
$ws.Cells.Item(2,5).Value = "❌"

# --- Row 3 (new) ---
$ws.Cells.Item(3,1).Value = 3
$ws.Cells.Item(3,2).Value = "PruebaEditar"
$ws.Cells.Item(3,3).Value = "18-09-2024"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,4).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3,4).Value = 45591
$ws.Cells.Item(3,5).Value = "✅"

# --- Row 4 (new) ---
$ws.Cells.Item(4,1).Value = 4
$ws.Cells.Item(4,2).Value = "Prueba10"
$ws.Cells.Item(4,3).Value = "18-09-2024"
$ws.Cells.Item(4,4).Value = "30-04-2025"
$ws.Cells.Item(4,5).Value = "❌"

# --- Row 5 (new) ---
$ws.Cells.Item(5,1).Value = 5
$ws.Cells.Item(5,2).Value = "Pruebaaa"
$ws.Cells.Item(5,3).Value = "22-09-2024"
$ws.Cells.Item(5,4).Value = "30-11-2024"
$ws.Cells.Item(5,5).Value = "❌"

# --- Column-level formatting (applied after data so it lands on every populated cell) ---
$ws.Columns.Item(1).NumberFormat = "#,##0"
$ws.Columns.Item(2).HorizontalAlignment = 1
$ws.Columns.Item(3).HorizontalAlignment = 1
$ws.Columns.Item(4).HorizontalAlignment = 1
$ws.Columns.Item(5).HorizontalAlignment = 1
